$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to stay plain text (this sheet stores
# every data column - other than the rank column A - as text, e.g. "213.16" or
# "26.301.43", and Excel would otherwise auto-coerce plain-decimal-looking
# strings into numbers). ClearFormats() afterwards removes the temporary "@"
# number format so the cell keeps its original (unstyled) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Apply the updated cryptocurrency price/volume figures (and the handful of
# ranking reshuffles) published by the "Updated cryptos list" GitHub Actions job.
Set-TextValue $ws.Range('D2') '26.301.43'
Set-TextValue $ws.Range('E2') '  +1.21%  '
Set-TextValue $ws.Range('D3') '1.627.08'
Set-TextValue $ws.Range('E3') '  +1.57%  '
Set-TextValue $ws.Range('E4') '  +0.11%  '
Set-TextValue $ws.Range('D5') '213.16'
Set-TextValue $ws.Range('E5') '  +0.75%  '
Set-TextValue $ws.Range('E6') '  +0.06%  '
Set-TextValue $ws.Range('E7') '  +0.43%  '
Set-TextValue $ws.Range('E8') '  +0.94%  '
Set-TextValue $ws.Range('D9') '0.0618'
Set-TextValue $ws.Range('E9') '  +0.69%  '
Set-TextValue $ws.Range('D10') '19.20'
Set-TextValue $ws.Range('E10') '  +5.92%  '
Set-TextValue $ws.Range('D11') '0.0821'
Set-TextValue $ws.Range('E11') '  +1.12%  '
Set-TextValue $ws.Range('D12') '1.854.30'
Set-TextValue $ws.Range('E12') '  +1.60%  '
Set-TextValue $ws.Range('D13') '1.624.50'
Set-TextValue $ws.Range('E13') '  +1.35%  '
Set-TextValue $ws.Range('D14') '4.05'
Set-TextValue $ws.Range('E14') '  +1.12%  '
Set-TextValue $ws.Range('D15') '0.521'
Set-TextValue $ws.Range('E15') '  +2.07%  '
Set-TextValue $ws.Range('D16') '26.321.70'
Set-TextValue $ws.Range('E16') '  +1.29%  '
Set-TextValue $ws.Range('D17') '62.65'
Set-TextValue $ws.Range('E17') '  +4.01%  '
Set-TextValue $ws.Range('D18') '0.0₃0732'
Set-TextValue $ws.Range('E18') '  +1.37%  '
Set-TextValue $ws.Range('E19') '  +0.17%  '
Set-TextValue $ws.Range('D20') '204.41'
Set-TextValue $ws.Range('E20') '  +1.26%  '
Set-TextValue $ws.Range('D21') '4.31'
Set-TextValue $ws.Range('E21') '  +1.91%  '
Set-TextValue $ws.Range('D22') '9.42'
Set-TextValue $ws.Range('E22') '  +1.74%  '
Set-TextValue $ws.Range('E23') '  +1.19%  '
Set-TextValue $ws.Range('D24') '1.94'
Set-TextValue $ws.Range('E24') '  +7.02%  '
Set-TextValue $ws.Range('D25') '144.09'
Set-TextValue $ws.Range('E25') '  +1.93%  '
Set-TextValue $ws.Range('E26') '  +0.10%  '
Set-TextValue $ws.Range('E27') '  +0.52%  '
Set-TextValue $ws.Range('D28') '15.39'
Set-TextValue $ws.Range('E28') '  +1.58%  '
Set-TextValue $ws.Range('E29') '  +2.54%  '
Set-TextValue $ws.Range('D30') '0.0524'
Set-TextValue $ws.Range('E30') '  +10.52%  '
Set-TextValue $ws.Range('E31') '  +0.60%  '
Set-TextValue $ws.Range('D32') '3.21'
Set-TextValue $ws.Range('E32') '  +3.12%  '
Set-TextValue $ws.Range('D33') '2.98'
Set-TextValue $ws.Range('E33') '  +0.55%  '
Set-TextValue $ws.Range('D34') '1.51'
Set-TextValue $ws.Range('E34') '  +2.69%  '
Set-TextValue $ws.Range('D35') '2.38'
Set-TextValue $ws.Range('E35') '  +1.41%  '
Set-TextValue $ws.Range('D36') '1.171.23'
Set-TextValue $ws.Range('E36') '  +3.86%  '
Set-TextValue $ws.Range('D37') '0.0167'
Set-TextValue $ws.Range('E37') '  +0.90%  '
Set-TextValue $ws.Range('B38') 'ARBITRUM'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D38') '0.810'
Set-TextValue $ws.Range('E38') '  +2.15%  '
Set-TextValue $ws.Range('B39') 'PaxDollar'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D39') '1.00'
Set-TextValue $ws.Range('E39') '  +0.15%  '
Set-TextValue $ws.Range('E40') '  +0.66%  '
Set-TextValue $ws.Range('D41') '0.502'
Set-TextValue $ws.Range('E41') '  +2.13%  '
Set-TextValue $ws.Range('D42') '0.794'
Set-TextValue $ws.Range('E42') '  +1.69%  '
Set-TextValue $ws.Range('D43') '5.34'
Set-TextValue $ws.Range('E43') '  +3.47%  '
Set-TextValue $ws.Range('D44') '1.767.46'
Set-TextValue $ws.Range('E44') '  +1.72%  '
Set-TextValue $ws.Range('D45') '93.21'
Set-TextValue $ws.Range('E45') '  +0.38%  '
Set-TextValue $ws.Range('B46') 'BabyDogeCoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.0₆0105'
Set-TextValue $ws.Range('E46') '  +13.98%  '
Set-TextValue $ws.Range('B47') 'RenderToken'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D47') '1.54'
Set-TextValue $ws.Range('E47') '  +1.40%  '
Set-TextValue $ws.Range('B48') 'Aave'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D48') '54.43'
Set-TextValue $ws.Range('E48') '  +1.83%  '
Set-TextValue $ws.Range('D49') '0.0510'
Set-TextValue $ws.Range('E49') '  +1.27%  '
Set-TextValue $ws.Range('D50') '0.410'
Set-TextValue $ws.Range('E50') '  +0.66%  '
Set-TextValue $ws.Range('E51') '  -0.01%  '
